$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(43).Insert()

$ws.Cells.Item(43, 1).Value = 6
$ws.Cells.Item(43, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(43, 3).Value = "Metropolitana"
$ws.Cells.Item(43, 4).Value = (Get-Date -Year 2021 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(43, 5).Value = 13
$ws.Cells.Item(43, 6).Value = 100112029
$ws.Cells.Item(43, 7).Value = "Orégano"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 32
$ws.Cells.Item(43, 11).Value = 9000
$ws.Cells.Item(43, 12).Value = 10000
$ws.Cells.Item(43, 13).Value = 9438
$ws.Cells.Item(43, 14).Value = "`$/docena de atados"
$ws.Cells.Item(43, 15).Value = "Región Metropolitana"
$ws.Cells.Item(43, 16).Value = 3146
$ws.Cells.Item(43, 17).Value = 3
$ws.Cells.Item(43, 18).Value = "Hortaliza"

Write-Host "done"
